$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 292.391276
$ws.Range("H2").Value = 877.173828
$ws.Range("I2").Value = 0.4546722242912879
$ws.Range("J2").Value = 0.4546722242912878
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.107333666666666
$ws.Range("N2").Value = 21.322001
$ws.Range("O2").Value = 0.7373665550576455
$ws.Range("P2").Value = 0.7373665550576454
$ws.Range("Q2").Value = 2078.122359754425
$ws.Range("R2").Value = 18703.10123778983
$ws.Range("S2").Value = 0.3352600917060641
$ws.Range("T2").Value = 0.335260091706064
$ws.Range("G3").Value = 292.391276
$ws.Range("H3").Value = 877.173828
$ws.Range("I3").Value = 0.4546722242912879
$ws.Range("J3").Value = 0.4546722242912878
$ws.Range("O3").Value = 0.1688878844614928
$ws.Range("P3").Value = 0.1688878844614928
$ws.Range("Q3").Value = 475.9772281285693
$ws.Range("R3").Value = 4283.795053157123
$ws.Range("S3").Value = 0.07678863008395699
$ws.Range("T3").Value = 0.07678863008395695
$ws.Range("G4").Value = 292.391276
$ws.Range("H4").Value = 877.173828
$ws.Range("I4").Value = 0.4546722242912879
$ws.Range("J4").Value = 0.4546722242912878
$ws.Range("M4").Value = 0.8135026666666666
$ws.Range("N4").Value = 2.440508
$ws.Range("O4").Value = 0.08439869112428164
$ws.Range("P4").Value = 0.08439869112428162
$ws.Range("Q4").Value = 237.8610827360693
$ws.Range("R4").Value = 2140.749744624624
$ws.Range("S4").Value = 0.03837374062075051
$ws.Range("T4").Value = 0.0383737406207505
$ws.Range("G5").Value = 292.391276
$ws.Range("H5").Value = 877.173828
$ws.Range("I5").Value = 0.4546722242912879
$ws.Range("J5").Value = 0.4546722242912878
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09009266666666667
$ws.Range("N5").Value = 0.270278
$ws.Range("O5").Value = 0.009346869356580103
$ws.Range("P5").Value = 0.009346869356580103
$ws.Range("Q5").Value = 26.34230976490933
$ws.Range("R5").Value = 237.080787884184
$ws.Range("S5").Value = 0.004249761880516355
$ws.Range("T5").Value = 0.004249761880516354
$ws.Range("H6").Value = 678.246018
$ws.Range("I6").Value = 0.3515604499097856
$ws.Range("J6").Value = 0.3515604499097856
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.107333666666666
$ws.Range("N6").Value = 21.322001
$ws.Range("O6").Value = 0.7373665550576455
$ws.Range("P6").Value = 0.7373665550576454
$ws.Range("Q6").Value = 1606.840252671335
$ws.Range("R6").Value = 14461.56227404202
$ws.Range("S6").Value = 0.2592289178444946
$ws.Range("T6").Value = 0.2592289178444945
$ws.Range("H7").Value = 678.246018
$ws.Range("I7").Value = 0.3515604499097856
$ws.Range("J7").Value = 0.3515604499097856
$ws.Range("O7").Value = 0.1688878844614928
$ws.Range("P7").Value = 0.1688878844614928
$ws.Range("S7").Value = 0.05937430064559431
$ws.Range("T7").Value = 0.05937430064559429
$ws.Range("H8").Value = 678.246018
$ws.Range("I8").Value = 0.3515604499097856
$ws.Range("J8").Value = 0.3515604499097856
$ws.Range("M8").Value = 0.8135026666666666
$ws.Range("N8").Value = 2.440508
$ws.Range("O8").Value = 0.08439869112428164
$ws.Range("P8").Value = 0.08439869112428162
$ws.Range("Q8").Value = 183.9183147663493
$ws.Range("R8").Value = 1655.264832897144
$ws.Range("S8").Value = 0.02967124182344948
$ws.Range("T8").Value = 0.02967124182344948
$ws.Range("H9").Value = 678.246018
$ws.Range("I9").Value = 0.3515604499097856
$ws.Range("J9").Value = 0.3515604499097856
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09009266666666667
$ws.Range("N9").Value = 0.270278
$ws.Range("O9").Value = 0.009346869356580103
$ws.Range("P9").Value = 0.009346869356580103
$ws.Range("Q9").Value = 20.36833080588934
$ws.Range("R9").Value = 183.314977253004
$ws.Range("S9").Value = 0.003285989596247289
$ws.Range("T9").Value = 0.003285989596247289
$ws.Range("G10").Value = 124.299764
$ws.Range("H10").Value = 372.899292
$ws.Range("I10").Value = 0.193287744251173
$ws.Range("J10").Value = 0.193287744251173
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.107333666666666
$ws.Range("N10").Value = 21.322001
$ws.Range("O10").Value = 0.7373665550576455
$ws.Range("P10").Value = 0.7373665550576454
$ws.Range("Q10").Value = 883.4398974359212
$ws.Range("R10").Value = 7950.959076923292
$ws.Range("S10").Value = 0.1425239181133507
$ws.Range("T10").Value = 0.1425239181133507
$ws.Range("G11").Value = 124.299764
$ws.Range("H11").Value = 372.899292
$ws.Range("I11").Value = 0.193287744251173
$ws.Range("J11").Value = 0.193287744251173
$ws.Range("O11").Value = 0.1688878844614928
$ws.Range("P11").Value = 0.1688878844614928
$ws.Range("Q11").Value = 202.3448097875373
$ws.Range("R11").Value = 1821.103288087836
$ws.Range("S11").Value = 0.03264395821891469
$ws.Range("T11").Value = 0.03264395821891468
$ws.Range("G12").Value = 124.299764
$ws.Range("H12").Value = 372.899292
$ws.Range("I12").Value = 0.193287744251173
$ws.Range("J12").Value = 0.193287744251173
$ws.Range("M12").Value = 0.8135026666666666
$ws.Range("N12").Value = 2.440508
$ws.Range("O12").Value = 0.08439869112428164
$ws.Range("P12").Value = 0.08439869112428162
$ws.Range("Q12").Value = 101.1181894800373
$ws.Range("R12").Value = 910.063705320336
$ws.Range("S12").Value = 0.0163132326251639
$ws.Range("T12").Value = 0.01631323262516389
$ws.Range("G13").Value = 124.299764
$ws.Range("H13").Value = 372.899292
$ws.Range("I13").Value = 0.193287744251173
$ws.Range("J13").Value = 0.193287744251173
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09009266666666667
$ws.Range("N13").Value = 0.270278
$ws.Range("O13").Value = 0.009346869356580103
$ws.Range("P13").Value = 0.009346869356580103
$ws.Range("Q13").Value = 11.19849720479733
$ws.Range("R13").Value = 100.786474843176
$ws.Range("S13").Value = 0.001806635293743781
$ws.Range("T13").Value = 0.001806635293743781
$ws.Range("G14").Value = 0.30841
$ws.Range("H14").Value = 0.92523
$ws.Range("I14").Value = 0.0004795815477534156
$ws.Range("J14").Value = 0.0004795815477534155
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.107333666666666
$ws.Range("N14").Value = 21.322001
$ws.Range("O14").Value = 0.7373665550576455
$ws.Range("P14").Value = 0.7373665550576454
$ws.Range("Q14").Value = 2.191972776136667
$ws.Range("R14").Value = 19.72775498523
$ws.Range("S14").Value = 0.0003536273937361497
$ws.Range("T14").Value = 0.0003536273937361497
$ws.Range("G15").Value = 0.30841
$ws.Range("H15").Value = 0.92523
$ws.Range("I15").Value = 0.0004795815477534156
$ws.Range("J15").Value = 0.0004795815477534155
$ws.Range("O15").Value = 0.1688878844614928
$ws.Range("P15").Value = 0.1688878844614928
$ws.Range("Q15").Value = 0.5020537511766666
$ws.Range("R15").Value = 4.51848376059
$ws.Range("S15").Value = 0.00008099551302684276
$ws.Range("T15").Value = 0.00008099551302684272
$ws.Range("G16").Value = 0.30841
$ws.Range("H16").Value = 0.92523
$ws.Range("I16").Value = 0.0004795815477534156
$ws.Range("J16").Value = 0.0004795815477534155
$ws.Range("M16").Value = 0.8135026666666666
$ws.Range("N16").Value = 2.440508
$ws.Range("O16").Value = 0.08439869112428164
$ws.Range("P16").Value = 0.08439869112428162
$ws.Range("Q16").Value = 0.2508923574266667
$ws.Range("R16").Value = 2.25803121684
$ws.Range("S16").Value = 0.00004047605491774544
$ws.Range("T16").Value = 0.00004047605491774544
$ws.Range("G17").Value = 0.30841
$ws.Range("H17").Value = 0.92523
$ws.Range("I17").Value = 0.0004795815477534156
$ws.Range("J17").Value = 0.0004795815477534155
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09009266666666667
$ws.Range("N17").Value = 0.270278
$ws.Range("O17").Value = 0.009346869356580103
$ws.Range("P17").Value = 0.009346869356580103
$ws.Range("Q17").Value = 0.02778547932666667
$ws.Range("R17").Value = 0.25006931394
$ws.Range("S17").Value = 0.000004482586072677657
$ws.Range("T17").Value = 0.000004482586072677657
